$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.463.12"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +4.40%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'4.044.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +3.48%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'518.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.92%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'146.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +1.32%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.726"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +18.43%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.757"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +5.03%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +1.49%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.0000326"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -2.49%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'47.04"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +11.53%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'10.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +6.24%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'4.685.58"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +3.33%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'4.030.93"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +1.88%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'21.08"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +6.69%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'14.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +0.41%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'  -1.33%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  -1.90%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'72.288.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +4.21%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'441.90"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +2.75%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'104.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +17.43%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'3.61"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +6.59%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'14.61"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +2.76%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'4.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -0.85%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'11.46"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -1.22%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'11.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +4.26%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'37.93"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +4.25%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  +2.22%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'3.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +10.92%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'13.62"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +3.33%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'0.129"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +1.99%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'676.31"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -1.28%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'6.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +13.94%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'67.44"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -0.49%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'42.38"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +5.73%  "
$ws.Range("E36").ClearFormats()
$ws.Range("B37").Value = "'PEPE"
$ws.Range("B37").ClearFormats()
$ws.Range("C37").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C37").ClearFormats()
$ws.Range("D37").Value = "'0.0₃0862"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +1.75%  "
$ws.Range("E37").ClearFormats()
$ws.Range("B38").Value = "'TheGraph"
$ws.Range("B38").ClearFormats()
$ws.Range("C38").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C38").ClearFormats()
$ws.Range("D38").Value = "'0.430"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -3.29%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'3.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +11.90%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.151"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +0.67%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +0.01%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.0496"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +2.90%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.998"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -0.23%  "
$ws.Range("E43").ClearFormats()
$ws.Range("B44").Value = "'WEMIXToken"
$ws.Range("B44").ClearFormats()
$ws.Range("C44").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C44").ClearFormats()
$ws.Range("D44").Value = "'3.22"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +2.86%  "
$ws.Range("E44").ClearFormats()
$ws.Range("B45").Value = "'Stellar"
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = "'0.158"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +12.15%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'2.72"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -3.06%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'3.44"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +2.42%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'3.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +2.03%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'9.04"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +6.50%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'3.32"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +1.34%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'  +1.32%  "
$ws.Range("E51").ClearFormats()
